$wb = $excel.ActiveWorkbook

# --- Work on the "map" worksheet ---
$ws = $wb.Worksheets.Item("map")

# Fix an error in the mapping: B3 should read "1A1a_Electricity" (was mistakenly
# "1A1a_Electricity-autoproducer", duplicating row 2's value combo).
$ws.Range("B3").Value = "1A1a_Electricity"

# Insert a new (blank) row above the old row 28 ("1A2g_Ind-Comb-Construction"),
# shifting it and everything below down by one row.
$ws.Rows.Item(28).Insert()

# Reflect the saved selection state shown in the workbook: bottom-right pane
# of the frozen "map" sheet has B27 selected.
$ws.Range("B27").Select()

# Make "map" the active/selected sheet in the workbook.
$ws.Activate()
